# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets hold the same table of convention listings; the diff bumps the
# same set of rows/values on each sheet.

$wb = $excel.ActiveWorkbook

# row -> new F value (only rows whose count changed)
$updates = @{
    2  = 1173
    3  = 934
    4  = 289
    5  = 57
    6  = 1122
    8  = 2408
    9  = 7879
    10 = 939
    11 = 463
    12 = 404
    13 = 165
    14 = 438
    15 = 8
    16 = 167
    17 = 8110
    18 = 326
    19 = 1403
    20 = 162
    23 = 182
    24 = 337
    25 = 183
    27 = 25
    29 = 35
    30 = 431
    31 = 1167
    32 = 15
    33 = 58
    34 = 102
    35 = 69
    36 = 88
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
